$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 26
$ws.Cells.Item($row, 1).Value = 57
$ws.Cells.Item($row, 2).Value = "Update index.py"
$ws.Cells.Item($row, 3).Value = "riya-morankar"
$ws.Cells.Item($row, 4).Value = "N/A"
$ws.Cells.Item($row, 5).Value = "edit1 to main"

$dateCell = $ws.Cells.Item($row, 6)
$dateCell.NumberFormat = "@"
$dateCell.Value = "2025-06-23"
$dateCell.Style = "Normal"
